$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 38 ---
# Add K38 = "Kagwe"
$ws.Range("K38").Value = "Kagwe"

# Swap O38/P38 values: O38 6->5, P38 5->6
$ws.Range("O38").Value = 5
$ws.Range("P38").Value = 6

# --- Add new row 39 ---
# Match A39's style/number-format to A38 (date, centered) instead of
# letting a custom NumberFormat string create a brand-new style entry.
$ws.Range("A39").NumberFormat = $ws.Range("A38").NumberFormat
$ws.Range("A39").Value = 43942
$ws.Range("B39").Value = 15
$ws.Range("C39").Value = 545
$ws.Range("D39").Value = "None"
$ws.Range("E39").Value = "Mombasa(7),Nairobi(6),Mandera(2)"
$ws.Range("F39").Value = 296
$ws.Range("G39").Value = "Community(15)"
$ws.Range("H39").Value = 5
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = "Mercy"
$ws.Range("L39").Value = "19-75"
$ws.Range("O39").Value = 8
$ws.Range("P39").Value = 7

# View state updates
$ws.Application.ActiveWindow.ScrollRow = 18
$ws.Range("O35").Select()
